# Refresh the cryptos table: Price (col D) and Volume(1h) (col E) for rows 2-51.
# A handful of Price cells are plain decimals (e.g. "608.58"); Range.Value on this
# host - like real Excel COM - auto-converts those to numbers unless the text is
# forced with a leading apostrophe, so those get the quote-prefix treatment to stay text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.698.60"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "2.692.26"
$ws.Range("E3").Value = "  +1.92%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "`'608.58"
$ws.Range("E5").Value = "  +2.20%  "
$ws.Range("D6").Value = "`'157.75"
$ws.Range("E6").Value = "  +0.94%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "`'0.588"
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("E9").Value = "  +5.21%  "
$ws.Range("D10").Value = "`'6.10"
$ws.Range("E10").Value = "  +5.61%  "
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("D13").Value = "`'30.16"
$ws.Range("E13").Value = "  +4.11%  "
$ws.Range("E14").Value = "  +10.57%  "
$ws.Range("D15").Value = "3.175.30"
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").Value = "65.507.40"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").Value = "2.695.76"
$ws.Range("E17").Value = "  +4.04%  "
$ws.Range("D18").Value = "`'12.72"
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("D20").Value = "`'359.81"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("E21").Value = "  +3.29%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("E23").Value = "  +2.78%  "
$ws.Range("D24").Value = "`'9.80"
$ws.Range("E24").Value = "  +2.87%  "
$ws.Range("D25").Value = "`'0.0000106"
$ws.Range("E25").Value = "  +12.61%  "
$ws.Range("D26").Value = "`'1.70"
$ws.Range("E26").Value = "  +3.22%  "
$ws.Range("E27").Value = "  -4.19%  "
$ws.Range("D28").Value = "`'0.170"
$ws.Range("E28").Value = "  +3.89%  "
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("D30").Value = "`'2.20"
$ws.Range("E30").Value = "  +4.93%  "
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").Value = "`'533.85"
$ws.Range("E32").Value = "  +2.49%  "
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("D34").Value = "`'6.72"
$ws.Range("E34").Value = "  +6.76%  "
$ws.Range("E35").Value = "  -2.88%  "
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("E37").Value = "  +2.43%  "
$ws.Range("D38").Value = "`'162.58"
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("E39").Value = "  -1.33%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "`'168.64"
$ws.Range("E42").Value = "  +2.15%  "
$ws.Range("D43").Value = "`'42.62"
$ws.Range("E43").Value = "  +0.96%  "
$ws.Range("E44").Value = "  +2.14%  "
$ws.Range("D45").Value = "`'0.0616"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "`'23.48"
$ws.Range("E46").Value = "  +2.25%  "
$ws.Range("D47").Value = "`'2.26"
$ws.Range("E47").Value = "  +2.21%  "
$ws.Range("E48").Value = "  +4.68%  "
$ws.Range("D49").Value = "`'0.657"
$ws.Range("E49").Value = "  +1.54%  "
$ws.Range("D50").Value = "`'20.98"
$ws.Range("E50").Value = "  +7.82%  "
$ws.Range("D51").Value = "`'0.0982"
$ws.Range("E51").Value = "  -0.37%  "
